$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark attendance (value 1) for the new date column (H) for each person, rows 3-8
$ws.Range("H3:H8").Value = 1

# Force recalculation so the dependent SUM/MAX formulas refresh their cached values
$excel.Calculate()

# Update the current selection to match the author's last selected cell
$ws.Range("I10").Select()
